# FDB-AOA Weights and Comparison results.
# Adds a new "dimension reduction" worksheet after "Scores" summarising the
# original data shape vs. the PSO- and FDB-AOA-reduced feature sets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the LAST sheet (after "Scores")
# ---------------------------------------------------------------------
$scores = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "dimension reduction"

# ---------------------------------------------------------------------
# 2. Column width / row height
# ---------------------------------------------------------------------
# Target stored width (Excel "characters" units, MDW=7) is 34.7109375;
# this runtime's ColumnWidth setter adds 5/6 internally, so back it out.
$ws.Columns.Item(1).ColumnWidth = 34.7109375 - (5/6)
$ws.Rows.Item(1).RowHeight = 45.75

# ---------------------------------------------------------------------
# 3. Header row (row 1)
#    Style creation order mirrors the authored workbook: wrap-only style
#    first, then bold-only, then bold+wrap, so the resulting cellXfs
#    index order lines up with the target file.
# ---------------------------------------------------------------------
$ws.Range("F1").WrapText = $true

$ws.Range("A1").Font.Bold = $true

$ws.Range("B1").Value = "Train"
$ws.Range("B1").Font.Bold = $true

$ws.Range("C1").Value = "Test"
$ws.Range("C1").Font.Bold = $true

$ws.Range("D1").Value = "Extraxted Features"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").WrapText = $true

$ws.Range("E1").Value = "Extracted Feature Rate"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").WrapText = $true

# ---------------------------------------------------------------------
# 4. Row 2 - Original Data Shape
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Original Data Shape"
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Value = "(507, 147)"
$ws.Range("C2").Value = "(168, 147) "
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("E2").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 5. Row 3 - Dimension Reduction (w/ PSO)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Dimension Reduction (w/ PSO)"
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Value = "(507, 130)"
$ws.Range("C3").Value = "(168, 130) "
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 0.12
$ws.Range("E3").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 6. Row 4 - Dimension Reduction (w/ FDB-AOA)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Dimension Reduction (w/ FDB-AOA)"
$ws.Range("A4").Font.Bold = $true
$ws.Range("B4").Value = "(507, 14)"
$ws.Range("C4").Value = "(168, 14) "
$ws.Range("D4").Formula = "=147-14"
$ws.Range("E4").Value = 0.9
$ws.Range("E4").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 7. Selection on the new sheet, then hand focus back to "Scores"
# ---------------------------------------------------------------------
$ws.Range("E6").Select() | Out-Null

[void]$scores.Activate()
$scores.Columns.Item(7).ColumnWidth = 32.85546875 - (5/6)
$scores.Columns.Item(8).ColumnWidth = 12.42578125 - (5/6)
$scores.Columns.Item(10).ColumnWidth = 19.28515625 - (5/6)
$scores.Range("G1:L4").Select() | Out-Null

Write-Host "dimension reduction sheet added"
